# front-end import task data connection
# Re-import task data: refresh the China value on row 2 and append a new
# row for the CB2027 / America container-boat record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-affirm the Country value on row 2 (re-synced from the data import).
$ws.Range("C2").Value = "China"

# New row of imported data.
$ws.Range("A5").Value = "CB2027"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "America"

# Carry the date/time formatting from the row above down into the new
# ArrivalTime / DepartureTime cells before writing their values.
$ws.Range("D4:E4").Copy()
$ws.Range("D5:E5").PasteSpecial(-4122)
$ws.Range("D5").Value = 45351.75
$ws.Range("E5").Value = 45352.75

$ws.Range("E5").Select() | Out-Null
